$d = $word.ActiveDocument

# Locate the whole paragraph describing the campaign (from "Συμμετέχετε" to
# the closing "...νυχτερινό ουρανό.") using a wildcard Find, then replace its
# entire contents with the updated, merged text that names the constellation
# "Αστερισμός Πήγασος" instead of "αστερισμού του Περσεύς".
$range = $d.Content
$found = $range.Find.Execute("Συμμετέχετε*νυχτερινό ουρανό.", $true, $false, $true, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $newText = "Συμμετέχετε σε μία παγκόσμια καμπάνια για να παρατηρήσετε και να καταγράψετε τη φωτεινότητα των πιο αμυδρά ορατών άστρων σαν μέσο για την μέτρηση της Φωτορρύπανσης σε μία δεδομένη περιοχή. Με τον εντοπισμό και την παρατήρηση του  Αστερισμός Πήγασος στον νυχτερινό ουρανό καθώς και με την σύγκριση των ανωτέρω με τα διαγράμματα για τα μεγέθη των άστρων,  άνθρωποι από όλον τον κόσμο θα μάθουν πώς τα φώτα στην κοινότητά τους συμβάλλουν στην Φωτορρύπανση. Με την κατάθεση των πορισμάτων τους στην ιστοσελίδα θα δημιουργηθεί ένα αρχείο σχετικά με το τι μπορεί να δει κανείς στον νυχτερινό ουρανό."

    $insertStart = $range.Start
    $range.Delete()
    $collapsed = $d.Range($insertStart, $insertStart)
    $collapsed.InsertAfter($newText)
}
